{"js": "// Replace the date line and all 100 answer cells in the practice-sheet\n// table. The document body, when flattened via `paragraphs`, yields the\n// title paragraph followed by the 100 table-cell paragraphs in the exact\n// same left-to-right, top-to-bottom order they appear in the OOXML, so we\n// can walk that flat list and swap each paragraph's text for its new\n// value (matched up by position, since a couple of the old answers repeat\n// verbatim but map to different new values depending on where they sit).\n\nconst oldValues = [\"2024-08-06 Tuesday\", \"45-3=42\", \"50+23=73\", \"52-22=30\", \"2-0=2\", \"16-5=11\", \"31-9=22\", \"58+28=86\", \"72-0=72\", \"90-61=29\", \"78-53=25\", \"46+32=78\", \"82-60=22\", \"99-63=36\", \"13+85=98\", \"78-42=36\", \"55+22=77\", \"38-13=25\", \"80+16=96\", \"4+70=74\", \"10+60=70\", \"64-42=22\", \"65-12=53\", \"67+27=94\", \"43-13=30\", \"10+64=74\", \"86+11=97\", \"0+41=41\", \"50-26=24\", \"90+5=95\", \"98-87=11\", \"35+39=74\", \"81+16=97\", \"81-41=40\", \"63-2=61\", \"94-8=86\", \"32-18=14\", \"62-22=40\", \"37-19=18\", \"36+25=61\", \"71-27=44\", \"62+5=67\", \"68-17=51\", \"83+1=84\", \"28-5=23\", \"27+10=37\", \"29+16=45\", \"40-21=19\", \"38-25=13\", \"98-40=58\", \"58-24=34\", \"25-22=3\", \"73-18=55\", \"69-2=67\", \"60+25=85\", \"78-55=23\", \"36+56=92\", \"49-30=19\", \"28-5=23\", \"1+13=14\", \"44+40=84\", \"96-16=80\", \"23+74=97\", \"66+21=87\", \"27-19=8\", \"66-18=48\", \"69-23=46\", \"90-72=18\", \"33-32=1\", \"14+10=24\", \"85-70=15\", \"94+2=96\", \"21+58=79\", \"59+2=61\", \"30+38=68\", \"35+13=48\", \"67+28=95\", \"41-4=37\", \"16+35=51\", \"28-5=23\", \"73-12=61\", \"54-34=20\", \"15+51=66\", \"55-27=28\", \"82-36=46\", \"26+33=59\", \"37+2=39\", \"51-0=51\", \"95-33=62\", \"35-28=7\", \"83-33=50\", \"43-25=18\", \"96-4=92\", \"44+20=64\", \"45+21=66\", \"48+42=90\", \"40-8=32\", \"53-32=21\", \"68-55=13\", \"45+27=72\", \"30+34=64\"];\nconst newValues = [\"2024-08-07 Wednesday\", \"20+61=81\", \"69-31=38\", \"99-59=40\", \"12+75=87\", \"27+41=68\", \"83+5=88\", \"98-70=28\", \"50-10=40\", \"51+27=78\", \"36+3=39\", \"67-2=65\", \"34+43=77\", \"21+77=98\", \"71-71=0\", \"26-1=25\", \"68+9=77\", \"26+0=26\", \"66+1=67\", \"48-24=24\", \"50+1=51\", \"34+31=65\", \"26+42=68\", \"60-56=4\", \"53+45=98\", \"81-2=79\", \"92-59=33\", \"29+31=60\", \"37+3=40\", \"18+52=70\", \"65-56=9\", \"15+59=74\", \"97-39=58\", \"86-42=44\", \"18+48=66\", \"60+4=64\", \"41+45=86\", \"26+7=33\", \"23+43=66\", \"67-2=65\", \"30+64=94\", \"25-20=5\", \"36+15=51\", \"72+7=79\", \"91-79=12\", \"24+6=30\", \"9+16=25\", \"55+6=61\", \"45+10=55\", \"21-14=7\", \"19+61=80\", \"82+13=95\", \"95-83=12\", \"26+19=45\", \"12+77=89\", \"29+17=46\", \"55-22=33\", \"19+13=32\", \"90-59=31\", \"11+35=46\", \"32+56=88\", \"20+42=62\", \"16+68=84\", \"29+41=70\", \"14+30=44\", \"32+61=93\", \"79+14=93\", \"11+16=27\", \"87-25=62\", \"95-4=91\", \"44-37=7\", \"46+12=58\", \"89-8=81\", \"87-50=37\", \"16+54=70\", \"82-77=5\", \"61+10=71\", \"72-47=25\", \"85-17=68\", \"51+32=83\", \"26-0=26\", \"32-22=10\", \"93-56=37\", \"13+71=84\", \"68-8=60\", \"54-37=17\", \"49-0=49\", \"83-82=1\", \"7+16=23\", \"68-56=12\", \"71-47=24\", \"76-2=74\", \"73+13=86\", \"30+53=83\", \"30+6=36\", \"90+3=93\", \"5+32=37\", \"80-30=50\", \"66-23=43\", \"96-3=93\", \"22+0=22\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== oldValues.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + oldValues.length + \" got \" + items.length\n  );\n}\n\n// Read current text for every paragraph first so we can confirm we are\n// about to edit the paragraph we think we are before touching anything.\nfor (let i = 0; i < items.length; i++) {\n  items[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < items.length; i++) {\n  const current = items[i].text.trim();\n  // Skip paragraphs that already hold the replacement value so a re-run of\n  // the script stays idempotent; everything else (expected to be the\n  // original value) gets replaced by position.\n  if (current !== newValues[i]) {\n    items[i].insertText(newValues[i], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and all 100 answer cells in the practice-sheet\n# table. Word COM flattens the whole document (including table cells) into\n# one Paragraphs collection in document order, but each table row also\n# contributes a trailing \"cell mark\" paragraph (its Range.Text is just the\n# paragraph mark + cell-end mark, i.e. empty once control characters are\n# stripped) after every 5 answer cells, so those markers are skipped while\n# walking the collection. The position (not the text) drives which new\n# value is written, since a few old answers repeat verbatim but land on\n# different replacements depending on where they sit in the sheet.\n\n$d = $word.ActiveDocument\n\n$oldValues = @(\n    \"2024-08-06 Tuesday\",\n    \"45-3=42\",\n    \"50+23=73\",\n    \"52-22=30\",\n    \"2-0=2\",\n    \"16-5=11\",\n    \"31-9=22\",\n    \"58+28=86\",\n    \"72-0=72\",\n    \"90-61=29\",\n    \"78-53=25\",\n    \"46+32=78\",\n    \"82-60=22\",\n    \"99-63=36\",\n    \"13+85=98\",\n    \"78-42=36\",\n    \"55+22=77\",\n    \"38-13=25\",\n    \"80+16=96\",\n    \"4+70=74\",\n    \"10+60=70\",\n    \"64-42=22\",\n    \"65-12=53\",\n    \"67+27=94\",\n    \"43-13=30\",\n    \"10+64=74\",\n    \"86+11=97\",\n    \"0+41=41\",\n    \"50-26=24\",\n    \"90+5=95\",\n    \"98-87=11\",\n    \"35+39=74\",\n    \"81+16=97\",\n    \"81-41=40\",\n    \"63-2=61\",\n    \"94-8=86\",\n    \"32-18=14\",\n    \"62-22=40\",\n    \"37-19=18\",\n    \"36+25=61\",\n    \"71-27=44\",\n    \"62+5=67\",\n    \"68-17=51\",\n    \"83+1=84\",\n    \"28-5=23\",\n    \"27+10=37\",\n    \"29+16=45\",\n    \"40-21=19\",\n    \"38-25=13\",\n    \"98-40=58\",\n    \"58-24=34\",\n    \"25-22=3\",\n    \"73-18=55\",\n    \"69-2=67\",\n    \"60+25=85\",\n    \"78-55=23\",\n    \"36+56=92\",\n    \"49-30=19\",\n    \"28-5=23\",\n    \"1+13=14\",\n    \"44+40=84\",\n    \"96-16=80\",\n    \"23+74=97\",\n    \"66+21=87\",\n    \"27-19=8\",\n    \"66-18=48\",\n    \"69-23=46\",\n    \"90-72=18\",\n    \"33-32=1\",\n    \"14+10=24\",\n    \"85-70=15\",\n    \"94+2=96\",\n    \"21+58=79\",\n    \"59+2=61\",\n    \"30+38=68\",\n    \"35+13=48\",\n    \"67+28=95\",\n    \"41-4=37\",\n    \"16+35=51\",\n    \"28-5=23\",\n    \"73-12=61\",\n    \"54-34=20\",\n    \"15+51=66\",\n    \"55-27=28\",\n    \"82-36=46\",\n    \"26+33=59\",\n    \"37+2=39\",\n    \"51-0=51\",\n    \"95-33=62\",\n    \"35-28=7\",\n    \"83-33=50\",\n    \"43-25=18\",\n    \"96-4=92\",\n    \"44+20=64\",\n    \"45+21=66\",\n    \"48+42=90\",\n    \"40-8=32\",\n    \"53-32=21\",\n    \"68-55=13\",\n    \"45+27=72\",\n    \"30+34=64\"\n)\n\n$newValues = @(\n    \"2024-08-07 Wednesday\",\n    \"20+61=81\",\n    \"69-31=38\",\n    \"99-59=40\",\n    \"12+75=87\",\n    \"27+41=68\",\n    \"83+5=88\",\n    \"98-70=28\",\n    \"50-10=40\",\n    \"51+27=78\",\n    \"36+3=39\",\n    \"67-2=65\",\n    \"34+43=77\",\n    \"21+77=98\",\n    \"71-71=0\",\n    \"26-1=25\",\n    \"68+9=77\",\n    \"26+0=26\",\n    \"66+1=67\",\n    \"48-24=24\",\n    \"50+1=51\",\n    \"34+31=65\",\n    \"26+42=68\",\n    \"60-56=4\",\n    \"53+45=98\",\n    \"81-2=79\",\n    \"92-59=33\",\n    \"29+31=60\",\n    \"37+3=40\",\n    \"18+52=70\",\n    \"65-56=9\",\n    \"15+59=74\",\n    \"97-39=58\",\n    \"86-42=44\",\n    \"18+48=66\",\n    \"60+4=64\",\n    \"41+45=86\",\n    \"26+7=33\",\n    \"23+43=66\",\n    \"67-2=65\",\n    \"30+64=94\",\n    \"25-20=5\",\n    \"36+15=51\",\n    \"72+7=79\",\n    \"91-79=12\",\n    \"24+6=30\",\n    \"9+16=25\",\n    \"55+6=61\",\n    \"45+10=55\",\n    \"21-14=7\",\n    \"19+61=80\",\n    \"82+13=95\",\n    \"95-83=12\",\n    \"26+19=45\",\n    \"12+77=89\",\n    \"29+17=46\",\n    \"55-22=33\",\n    \"19+13=32\",\n    \"90-59=31\",\n    \"11+35=46\",\n    \"32+56=88\",\n    \"20+42=62\",\n    \"16+68=84\",\n    \"29+41=70\",\n    \"14+30=44\",\n    \"32+61=93\",\n    \"79+14=93\",\n    \"11+16=27\",\n    \"87-25=62\",\n    \"95-4=91\",\n    \"44-37=7\",\n    \"46+12=58\",\n    \"89-8=81\",\n    \"87-50=37\",\n    \"16+54=70\",\n    \"82-77=5\",\n    \"61+10=71\",\n    \"72-47=25\",\n    \"85-17=68\",\n    \"51+32=83\",\n    \"26-0=26\",\n    \"32-22=10\",\n    \"93-56=37\",\n    \"13+71=84\",\n    \"68-8=60\",\n    \"54-37=17\",\n    \"49-0=49\",\n    \"83-82=1\",\n    \"7+16=23\",\n    \"68-56=12\",\n    \"71-47=24\",\n    \"76-2=74\",\n    \"73+13=86\",\n    \"30+53=83\",\n    \"30+6=36\",\n    \"90+3=93\",\n    \"5+32=37\",\n    \"80-30=50\",\n    \"66-23=43\",\n    \"96-3=93\",\n    \"22+0=22\"\n)\n\n$expectedCount = $oldValues.Count\n\n$index = 0\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $visible = $r.Text -replace \"[\\x07\\x0d\\x0c\\x02]\", \"\"\n    if ($visible -eq \"\") {\n        # Row-end / cell-mark-only paragraph - nothing to replace.\n        continue\n    }\n    if ($index -ge $expectedCount) {\n        break\n    }\n    # Only write when needed: leave an already-updated cell alone (so\n    # re-running this script is a no-op) and otherwise overwrite by\n    # position regardless of whether the text still matches the recorded\n    # original, which keeps us robust to tiny whitespace drift.\n    if ($visible -ne $newValues[$index]) {\n        $r.Text = $newValues[$index]\n    }\n    $index = $index + 1\n}\n\nif ($index -ne $expectedCount) {\n    throw \"Expected to update $expectedCount paragraphs but updated $index\"\n}\n"}
